# Actualización automática 2025-11-13 16:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D4").Value = 561.7
$ws1.Range("M4").Value = 1682.73
$ws1.Range("M12").Value = 1706.31
$ws1.Range("I31").Value = 28.8
$ws1.Range("I56").Value = "3 de 54"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 3345.36
$ws2.Range("F12").Value = 1706.31
$ws2.Range("F31").Value = 28.8
$ws2.Range("F60").Value = 20958.88

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 2885.86
$ws3.Range("E3").Value = 3737.4
$ws3.Range("F3").Value = 0.4357159465278428

$ws3.Range("D7").Value = 365.4
$ws3.Range("E7").Value = 954.6
$ws3.Range("F7").Value = 0.2768181818181818

$ws3.Range("D12").Value = 9132.51
$ws3.Range("E12").Value = 55811.49
$ws3.Range("F12").Value = 0.1406213045084996

$ws3.Range("D14").Value = 20935.48
$ws3.Range("E14").Value = 78020.77685923838
$ws3.Range("F14").Value = 0.2115629740298277
